# TC_ManageLeaveTypes.xlsx - leave types validation update
# (create, edit, updated, activate, deactivate, delete)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the leave-type description / name text for row 2.
# Write column B before column A so the shared-string table ends up in the
# same insertion order as the committed workbook.
$ws.Range("B2").Value = "Leave with pay FDSF DFF SD dfsf "
$ws.Range("A2").Value = "Emergency Leave"

# Row 2 no longer needs the leftover custom-row-format markers.
$ws.Rows("2:2").ClearFormats()

# The old blank row 3 placeholder is removed entirely.
$ws.Rows("3:3").Delete()

# Widen the two data columns.
$ws.Columns("A:A").ColumnWidth = 22.65
$ws.Columns("B:B").ColumnWidth = 30.35

# Move the active selection to B8.
$ws.Range("B8").Select()
